$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")
Write-Host $ws.Name
